$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(39, 6508, 45747.42708333334),
    @(40, 6478, 45747.4375),
    @(41, 6320, 45747.44791666666),
    @(42, 6226, 45747.45833333334),
    @(43, 6202, 45747.46875),
    @(44, 6204, 45747.47916666666),
    @(45, 6111, 45747.48958333334),
    @(46, 6112, 45747.5)
)

foreach ($row in $data) {
    $r = $row[0]
    $aVal = $row[1]
    $bVal = $row[2]

    $ws.Cells.Item($r, 1).Value = $aVal

    $bCell = $ws.Cells.Item($r, 2)
    $bCell.Value = $bVal
    $bCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
